$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: duplicate the current (old) row 24 data into the new row 25,
# before row 24 itself gets overwritten with this week's new data.
$ws.Range("A25").Value = 3
$ws.Range("B25").Value = "Femacal de La Calera"
$ws.Range("C25").Value = "Coquimbo"
$ws.Range("D25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D25").Value = 44312
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100104
$ws.Range("H25").Value = "Frutos de pepita"
$ws.Range("I25").Value = 100104001
$ws.Range("J25").Value = "Granada"
$ws.Range("K25").Value = "Wonderfull"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 68
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 14000
$ws.Range("P25").Value = 14000
$ws.Range("Q25").Value = "$/caja 14 kilos granel"
$ws.Range("R25").Value = "Provincia de Limarí"
$ws.Range("S25").Value = 1000
$ws.Range("T25").Value = 14

# Step 2: update row 24 with the new weekly price entry.
$ws.Range("D24").Value = 44627
$ws.Range("M24").Value = 56
$ws.Range("N24").Value = 17000
$ws.Range("O24").Value = 17000
$ws.Range("P24").Value = 17000
$ws.Range("Q24").Value = "$/caja 14 kilos empedrada"
$ws.Range("S24").Value = 1214
